$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the D column "=B-C" shared formula down through the new rows first
# (setting the formula before the B/C values keeps the shared-formula block
# intact with a proper master/ref cell once the referenced values land).
$ws.Range("D28:D31").Formula = "=B28-C28"

# --- Fill in the newly-reported daily rows (28-31) ---
# Row 28 (serial 45985)
$ws.Range("B28").Value = 2047
$ws.Range("C28").Value = 1915
$ws.Range("E28").Value = 52
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 352
$ws.Range("I28").Value = 5
$ws.Range("J28").Value = 0

# Row 29 (serial 45986)
$ws.Range("B29").Value = 1902
$ws.Range("C29").Value = 1747
$ws.Range("E29").Value = 46
$ws.Range("F29").Value = 8
$ws.Range("G29").Value = 8
$ws.Range("H29").Value = 320
$ws.Range("I29").Value = 12
$ws.Range("J29").Value = 32

# Row 30 (serial 45987)
$ws.Range("B30").Value = 2455
$ws.Range("C30").Value = 2257
$ws.Range("E30").Value = 49
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 352
$ws.Range("I30").Value = 8
$ws.Range("J30").Value = 683

# Row 31 (serial 45988)
$ws.Range("B31").Value = 2595
$ws.Range("C31").Value = 2352
$ws.Range("E31").Value = 54
$ws.Range("F31").Value = 5
$ws.Range("G31").Value = 8
$ws.Range("H31").Value = 345
$ws.Range("I31").Value = 11
$ws.Range("J31").Value = 66

# --- Update the active selection on the sheet (bottomRight pane) ---
$ws.Range("H9").Select()

# --- Shrink the saved window width ---
$excel.ActiveWindow.Width = 13545
